# Mark additional "Love Babbar 450" problems as done: set the status
# column (C) from "<->" to "yes" for the rows whose code/PDF solutions
# were added in this commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$doneCells = @(
    "C80", "C81", "C82", "C83", "C84",
    "C88", "C89", "C90", "C91", "C92", "C93", "C94", "C95", "C96",
    "C98",
    "C268",
    "C271", "C272",
    "C339",
    "C341", "C342",
    "C345", "C346", "C347", "C348", "C349",
    "C351", "C352", "C353",
    "C356",
    "C359", "C360"
)

foreach ($addr in $doneCells) {
    $ws.Range($addr).Value = "yes"
}

# Reflect the author's final on-screen position/selection when the
# workbook was saved.
$null = $ws.Range("B338").Select()

Write-Output "Updated $($doneCells.Count) cells to 'yes'"
